$wb = $excel.ActiveWorkbook

$wsConcepts = $wb.Worksheets.Item("Concepts")
$wsCriteria = $wb.Worksheets.Item("Criteria")
$wsMatrix   = $wb.Worksheets.Item("Trade-off Matrix")

# --- Concepts sheet: replace the 6 old concepts with 4 new ones ---
$wsConcepts.Range("B2").Value = "Polar Bear"
$wsConcepts.Range("B3").Value = "Grizzly Bear"
$wsConcepts.Range("B4").Value = "Beluga"
$wsConcepts.Range("B5").Value = "Dolphin"
$wsConcepts.Rows("6:7").Delete() | Out-Null

# --- Criteria sheet: new weights, and swap Risk/Cost order ---
$wsCriteria.Range("B2").Value = 0.3
$wsCriteria.Range("B3").Value = 0.3
$wsCriteria.Range("A4").Value = "Cost"
$wsCriteria.Range("B4").Value = 0.2
$wsCriteria.Range("A5").Value = "Risk"
$wsCriteria.Range("B5").Value = 0.2

# --- Trade-off Matrix sheet: new scores, two fewer concept rows ---
$wsMatrix.Range("B2").Value = 4
$wsMatrix.Range("C2").Value = 5
$wsMatrix.Range("D2").Value = 5
$wsMatrix.Range("E2").Value = 3

$wsMatrix.Range("B3").Value = 1
$wsMatrix.Range("C3").Value = 1
$wsMatrix.Range("D3").Value = 1
$wsMatrix.Range("E3").Value = 5

$wsMatrix.Range("B4").Value = 5
$wsMatrix.Range("C4").Value = 5
$wsMatrix.Range("D4").Value = 5
$wsMatrix.Range("E4").Value = 1

$wsMatrix.Range("B5").Value = 3
$wsMatrix.Range("C5").Value = 3
$wsMatrix.Range("D5").Value = 3
$wsMatrix.Range("E5").Value = 3

$wsMatrix.Rows("6:7").Delete() | Out-Null

# re-enter the dynamic-array formula in A2 so the spill shrinks to match
# the now-smaller Concepts list (A2:A5 instead of A2:A7)
$wsMatrix.Range("A2:A7").ClearContents() | Out-Null
$wsMatrix.Range("A2:A5").FormulaArray = "=Concepts!A2:A5"

$wb.Application.Calculate()

# --- restore selections / active sheet (Trade-off Matrix last = active) ---
$wsConcepts.Range("B6").Select() | Out-Null
$wsCriteria.Range("A6").Select() | Out-Null
$wsMatrix.Range("E2").Select() | Out-Null
